$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new Price values below are numeric-looking strings (e.g. "1.002").
# Format those D-column cells as Text first so Excel keeps them as
# strings (matching the source inlineStr cells) instead of silently
# converting them into real numbers.
$textRows = @(4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22, 23, 25, 26, 27, 28, 29, 30, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 44, 45, 46, 48, 49, 50)
foreach ($r in $textRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Range("D2").Value = "22.419.15"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "1.568.79"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "285.06"
$ws.Range("E6").Value = "  -2.26%  "
$ws.Range("D7").Value = "0.3641"
$ws.Range("E7").Value = "  -2.38%  "
$ws.Range("D8").Value = "48.65"
$ws.Range("E8").Value = "  -2.61%  "
$ws.Range("D9").Value = "0.3331"
$ws.Range("E9").Value = "  -2.09%  "
$ws.Range("D10").Value = "1.123"
$ws.Range("E10").Value = "  -1.89%  "
$ws.Range("D11").Value = "0.07390"
$ws.Range("E11").Value = "  -2.40%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "20.75"
$ws.Range("E13").Value = "  -2.69%  "
$ws.Range("D14").Value = "5.954"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").Value = "6.903"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").Value = "1.573.47"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "0.00001104"
$ws.Range("E17").Value = "  -1.61%  "
$ws.Range("D18").Value = "88.15"
$ws.Range("E18").Value = "  -3.10%  "
$ws.Range("D19").Value = "0.06701"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "6.343"
$ws.Range("E21").Value = "  +0.62%  "
$ws.Range("D22").Value = "16.18"
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("D23").Value = "12.00"
$ws.Range("E23").Value = "  -1.52%  "
$ws.Range("D24").Value = "22.408.31"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").Value = "2.383"
$ws.Range("E25").Value = "  +1.80%  "
$ws.Range("D26").Value = "2.535"
$ws.Range("E26").Value = "  -5.85%  "
$ws.Range("D27").Value = "150.28"
$ws.Range("E27").Value = "  +1.11%  "
$ws.Range("D28").Value = "19.38"
$ws.Range("E28").Value = "  -3.52%  "
$ws.Range("D29").Value = "5.006"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").Value = "123.88"
$ws.Range("E30").Value = "  -1.40%  "
$ws.Range("D31").Value = "1.756.83"
$ws.Range("E31").Value = "  +0.50%  "
$ws.Range("D32").Value = "1.044"
$ws.Range("E32").Value = "  -1.22%  "
$ws.Range("D33").Value = "6.103"
$ws.Range("E33").Value = "  -1.60%  "
$ws.Range("D34").Value = "1.990"
$ws.Range("E34").Value = "  +0.31%  "
$ws.Range("D35").Value = "9.803"
$ws.Range("E35").Value = "  -0.53%  "
$ws.Range("D36").Value = "0.08252"
$ws.Range("E36").Value = "  -1.52%  "
$ws.Range("D37").Value = "0.02414"
$ws.Range("E37").Value = "  -3.34%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "0.06430"
$ws.Range("E38").Value = "  -1.35%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "0.2234"
$ws.Range("E39").Value = "  -2.96%  "
$ws.Range("D40").Value = "5.355"
$ws.Range("E40").Value = "  -2.28%  "
$ws.Range("D41").Value = "1.289"
$ws.Range("E41").Value = "  -4.98%  "
$ws.Range("D42").Value = "11.17"
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").Value = "1.004"
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("D45").Value = "13.75"
$ws.Range("E45").Value = "  -1.74%  "
$ws.Range("D46").Value = "0.6021"
$ws.Range("E46").Value = "  +3.54%  "
$ws.Range("E47").Value = "  -1.60%  "
$ws.Range("D48").Value = "2.032"
$ws.Range("E48").Value = "  -1.96%  "
$ws.Range("D49").Value = "123.62"
$ws.Range("E49").Value = "  -4.64%  "
$ws.Range("D50").Value = "1.214"
$ws.Range("E50").Value = "  -0.93%  "
